$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "Metadata" sheet currently lists two "Contact" rows (10 and 11).
# The commit adds all IG authors as contact, i.e. it inserts two more
# "Contact" / "No display for ContactDetail" rows right after the
# existing ones (rows 12 and 13), pushing the remaining metadata rows
# (Jurisdiction, Description, Purpose, Copyright, Immutable) down by two
# rows (from 12-16 to 14-18).

# Shift the trailing rows down by two, starting from the bottom so we
# never overwrite data we still need to move. For each row we copy the
# value first, then copy the formatting (style only) from the source row
# so the existing shared style is reused instead of a new style being
# created.
for ($src = 16; $src -ge 12; $src--) {
    $dst = $src + 2

    $ws.Range("A$dst").Value = $ws.Range("A$src").Value()
    $ws.Range("B$dst").Value = $ws.Range("B$src").Value()

    $ws.Range("A${src}:B${src}").Copy()
    $ws.Range("A${dst}:B${dst}").PasteSpecial(-4122)
}

# Fill the two newly freed rows (12 and 13) with another "Contact" entry
# each, copying the look of the existing Contact row (row 11).
$ws.Range("A12").Value = $ws.Range("A11").Value()
$ws.Range("B12").Value = $ws.Range("B11").Value()
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)

$ws.Range("A13").Value = $ws.Range("A11").Value()
$ws.Range("B13").Value = $ws.Range("B11").Value()
$ws.Range("A11:B11").Copy()
$ws.Range("A13:B13").PasteSpecial(-4122)

$excel.CutCopyMode = $false
